$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells being updated so numeric-looking
# strings (e.g. "557.05", "108.00") are preserved verbatim as text, matching
# the source data which stores Price as inline strings (not numbers).
$dCells = @("D2","D3","D5","D6","D9","D10","D11","D13","D14","D16","D17","D18","D20","D21","D22","D23","D24","D25","D28","D31","D34","D36","D38","D39","D40","D45","D47","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '65.782.09'
$ws.Range("E2").Value = '  -4.22%  '

# Row 3
$ws.Range("D3").Value = '3.288.07'
$ws.Range("E3").Value = '  -5.62%  '

# Row 4
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
$ws.Range("D5").Value = '557.05'
$ws.Range("E5").Value = '  -4.02%  '

# Row 6
$ws.Range("D6").Value = '181.13'
$ws.Range("E6").Value = '  -4.26%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("E8").Value = '  -1.99%  '

# Row 9
$ws.Range("D9").Value = '3.289.32'
$ws.Range("E9").Value = '  -5.18%  '

# Row 10
$ws.Range("D10").Value = '0.186'
$ws.Range("E10").Value = '  -7.24%  '

# Row 11
$ws.Range("D11").Value = '0.589'
$ws.Range("E11").Value = '  -3.91%  '

# Row 12
$ws.Range("E12").Value = '  -7.33%  '

# Row 13
$ws.Range("D13").Value = '0.0000265'
$ws.Range("E13").Value = '  -6.31%  '

# Row 14
$ws.Range("D14").Value = '640.98'
$ws.Range("E14").Value = '  +1.14%  '

# Row 15
$ws.Range("E15").Value = '  -5.26%  '

# Row 16
$ws.Range("D16").Value = '3.825.05'
$ws.Range("E16").Value = '  -4.91%  '

# Row 17
$ws.Range("D17").Value = '65.858.15'
$ws.Range("E17").Value = '  -4.53%  '

# Row 18
$ws.Range("D18").Value = '17.93'
$ws.Range("E18").Value = '  -0.80%  '

# Row 19
$ws.Range("E19").Value = '  -3.11%  '

# Row 20
$ws.Range("D20").Value = '3.292.59'
$ws.Range("E20").Value = '  -5.41%  '

# Row 21
$ws.Range("D21").Value = '11.42'
$ws.Range("E21").Value = '  -7.52%  '

# Row 22
$ws.Range("D22").Value = '0.908'
$ws.Range("E22").Value = '  -3.68%  '

# Row 23
$ws.Range("D23").Value = '17.94'
$ws.Range("E23").Value = '  +1.24%  '

# Row 24
$ws.Range("D24").Value = '108.00'
$ws.Range("E24").Value = '  +8.49%  '

# Row 25
$ws.Range("D25").Value = '5.00'
$ws.Range("E25").Value = '  -6.90%  '

# Row 26
$ws.Range("E26").Value = '  -7.21%  '

# Row 27
$ws.Range("E27").Value = '  -6.01%  '

# Row 28
$ws.Range("D28").Value = '9.55'
$ws.Range("E28").Value = '  -4.76%  '

# Row 29
$ws.Range("E29").Value = '  -5.06%  '

# Row 30
$ws.Range("E30").Value = '  -6.12%  '

# Row 31
$ws.Range("D31").Value = '3.99'
$ws.Range("E31").Value = '  -1.40%  '

# Row 32
$ws.Range("E32").Value = '  -4.85%  '

# Row 33
$ws.Range("E33").Value = '  -3.92%  '

# Row 34
$ws.Range("D34").Value = '555.41'
$ws.Range("E34").Value = '  +10.55%  '

# Row 35
$ws.Range("E35").Value = '  -2.63%  '

# Row 36
$ws.Range("D36").Value = '57.25'
$ws.Range("E36").Value = '  -5.73%  '

# Row 37
$ws.Range("E37").Value = '  +0.08%  '

# Row 38
$ws.Range("D38").Value = '3.673.41'
$ws.Range("E38").Value = '  -0.87%  '

# Row 39
$ws.Range("B39").Value = 'CoreDAO'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D39").Value = '3.77'
$ws.Range("E39").Value = '  +40.78%  '

# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '3.51'
$ws.Range("E40").Value = '  -0.85%  '

# Row 41
$ws.Range("E41").Value = '  -6.19%  '

# Row 42
$ws.Range("E42").Value = '  -9.27%  '

# Row 43
$ws.Range("E43").Value = '  -3.71%  '

# Row 44
$ws.Range("E44").Value = '  -6.58%  '

# Row 45
$ws.Range("D45").Value = '32.11'

# Row 46
$ws.Range("E46").Value = '  -4.91%  '

# Row 47
$ws.Range("D47").Value = '3.25'
$ws.Range("E47").Value = '  -2.86%  '

# Row 48
$ws.Range("E48").Value = '  -3.13%  '

# Row 49
$ws.Range("D49").Value = '2.61'
$ws.Range("E49").Value = '  -6.22%  '

# Row 50
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.25%  '

# Row 51
$ws.Range("D51").Value = '7.64'
$ws.Range("E51").Value = '  -5.28%  '

# Restore default (unstyled) cell style on the Price cells now that the
# text value is committed, so the saved style table matches the original
# (unstyled) data cells.
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}